$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Softube Dynamics related labels (adds 3 new shared strings: "FX Zone", "Home", "Lock Item")
$ws.Range("I21").Value = "FX Zone"
$ws.Range("I22").Value = "Home"

# N40 gets a new bold label, matching the style used by the other filled
# "Rightcol" cells in column N (style index 5 = bold Consolas 9)
$ws.Range("N40").Value = "Lock Item"
$ws.Range("N40").Font.Bold = $true

# Column N was manually widened (no longer auto (best-fit) sized)
$ws.Columns.Item(14).ColumnWidth = 11.7

# Selection moved to N40 before the file was saved
$ws.Range("N40").Select()
